# Update gh-pages output data (refreshed counts/prices) across the
# 展览 (Exhibitions), 演出 (Shows) and 全部类型 (All types) sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 34      # 张家港·心动次元动漫游戏嘉年华          33 -> 34
$ws1.Range("F4").Value  = 16141   # 苏州·I COME ACG动漫品牌博览会         16135 -> 16141
$ws1.Range("F5").Value  = 423     # 苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场  422 -> 423
$ws1.Range("F8").Value  = 15547   # 苏州·理想乡动漫游戏展                 15542 -> 15547
$ws1.Range("G8").Value  = 76      #                                      19.9 -> 76
$ws1.Range("F10").Value = 9176    # 苏州·第四届-OCG国朝动漫游戏嘉年华      9171 -> 9176
$ws1.Range("F11").Value = 448     # 常熟·CDW·动漫展06                    447 -> 448
$ws1.Range("F14").Value = 116     # 苏州·OCG国潮动漫游戏嘉年华火只木南内场  115 -> 116
$ws1.Range("F17").Value = 215     # 苏州·第二届百合Only同人展交流          214 -> 215
$ws1.Range("F28").Value = 512     # 苏州·漫遇引力动漫游戏展               510 -> 512
$ws1.Range("F30").Value = 42      # 苏州·明日方舟同人only （聚会）         41 -> 42
$ws1.Range("F32").Value = 73      # 苏州·无限次元夜场                    72 -> 73
$ws1.Range("F37").Value = 469     # 苏州·绘时国乙2.0光夜同人only          467 -> 469
$ws1.Range("F39").Value = 5644    # 苏州·COME IN JOY 动漫品牌国潮文化节    5643 -> 5644

# ---------------------------------------------------------------
# Sheet "演出" (Shows)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 77      # 苏州·乐队番同人only live Band Set二次元乐队拼盘  76 -> 77

# ---------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 34      # 张家港·心动次元动漫游戏嘉年华          33 -> 34
$ws4.Range("F4").Value  = 16141   # 苏州·I COME ACG动漫品牌博览会         16136 -> 16141
$ws4.Range("F5").Value  = 423     # 苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场  422 -> 423
$ws4.Range("F8").Value  = 15547   # 苏州·理想乡动漫游戏展                 15542 -> 15547
$ws4.Range("G8").Value  = 76      #                                      19.9 -> 76
$ws4.Range("F10").Value = 9176    # 苏州·第四届-OCG国朝动漫游戏嘉年华      9171 -> 9176
$ws4.Range("F11").Value = 448     # 常熟·CDW·动漫展06                    447 -> 448
$ws4.Range("F14").Value = 116     # 苏州·OCG国潮动漫游戏嘉年华火只木南内场  115 -> 116
$ws4.Range("F17").Value = 215     # 苏州·第二届百合Only同人展交流          214 -> 215
$ws4.Range("F28").Value = 512     # 苏州·漫遇引力动漫游戏展               510 -> 512
$ws4.Range("F30").Value = 42      # 苏州·明日方舟同人only （聚会）         41 -> 42
$ws4.Range("F31").Value = 77      # 苏州·乐队番同人only live Band Set二次元乐队拼盘  76 -> 77
$ws4.Range("F34").Value = 73      # 苏州·无限次元夜场                    72 -> 73
$ws4.Range("F39").Value = 469     # 苏州·绘时国乙2.0光夜同人only          467 -> 469
$ws4.Range("F41").Value = 5644    # 苏州·COME IN JOY 动漫品牌国潮文化节    5643 -> 5644
